$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.81"
$ws.Range("E2").Value = "'0.25%"
$ws.Range("D3").Value = "'36.30"
$ws.Range("E3").Value = "'-0.89%"
$ws.Range("D4").Value = "'5.063"
$ws.Range("E4").Value = "'0.84%"
$ws.Range("D5").Value = "'0.07892"
$ws.Range("E5").Value = "'0.63%"
$ws.Range("D6").Value = "'2.130"
$ws.Range("E6").Value = "'-1.63%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.987"
$ws.Range("E7").Value = "'-0.51%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9283"
$ws.Range("E8").Value = "'0.71%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09735"
$ws.Range("E9").Value = "'-2.25%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1861"
$ws.Range("E10").Value = "'-0.50%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09026"
$ws.Range("E11").Value = "'4.19%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03769"
$ws.Range("E12").Value = "'4.19%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09894"
$ws.Range("E13").Value = "'-0.40%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001437"
$ws.Range("E14").Value = "'-2.34%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005652"
$ws.Range("E15").Value = "'-0.64%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.467"
$ws.Range("E16").Value = "'0.16%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.157"
$ws.Range("E17").Value = "'2.42%"
$ws.Range("E18").Value = "'14.03%"
$ws.Range("E19").Value = "'-0.77%"
$ws.Range("D20").Value = "'0.1313"
$ws.Range("E20").Value = "'-2.47%"
$ws.Range("E21").Value = "'4.04%"
$ws.Range("E22").Value = "'2.47%"
$ws.Range("D23").Value = "'0.04595"
$ws.Range("E23").Value = "'-0.42%"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'0.60%"
$ws.Range("E25").Value = "'-7.50%"
$ws.Range("D26").Value = "'0.0001306"
$ws.Range("E26").Value = "'-6.88%"
$ws.Range("E27").Value = "'74.13%"
$ws.Range("D39").Value = "'0.01958"
$ws.Range("E39").Value = "'8.31%"
$ws.Range("D40").Value = "'0.04966"
$ws.Range("E40").Value = "'4.86%"
$ws.Range("D41").Value = "'0.007828"
$ws.Range("E41").Value = "'-0.43%"
$ws.Range("D42").Value = "'0.1393"
$ws.Range("E42").Value = "'-0.86%"
$ws.Range("D43").Value = "'0.007831"
$ws.Range("E43").Value = "'3.03%"
$ws.Range("D44").Value = "'0.002134"
$ws.Range("E44").Value = "'-2.25%"
$ws.Range("D45").Value = "'0.01125"
$ws.Range("E45").Value = "'7.68%"
$ws.Range("D46").Value = "'0.00006278"
$ws.Range("E46").Value = "'-1.04%"
$ws.Range("E47").Value = "'0.28%"
$ws.Range("E48").Value = "'0.12%"
$ws.Range("D49").Value = "'51.69"
$ws.Range("E49").Value = "'43.00%"
$ws.Range("D50").Value = "'0.001906"
$ws.Range("E50").Value = "'-29.23%"
$ws.Range("E51").Value = "'0.28%"